$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = 111396060
$ws.Range("I2").NumberFormat = "@"
$ws.Range("I2").Value = "90"
$ws.Range("I2").Style = "Normal"
$ws.Range("P2").Value = "S om järnvägen - 3, Vg"
$ws.Range("Q2").Value = 432076.641898193
$ws.Range("R2").Value = 6419661.774153749

# Row 3
$ws.Range("A3").Value = 111396045
$ws.Range("I3").NumberFormat = "@"
$ws.Range("I3").Value = "70"
$ws.Range("I3").Style = "Normal"
$ws.Range("J3").Value = "stjälkar/strån/skott"
$ws.Range("P3").Value = "S om järnvägen, Vg"
$ws.Range("Q3").Value = 431889.3909100805
$ws.Range("R3").Value = 6419670.266848063

# Row 4
$ws.Range("A4").Value = 111396053
$ws.Range("I4").NumberFormat = "@"
$ws.Range("I4").Value = "80"
$ws.Range("I4").Style = "Normal"
$ws.Range("J4").NumberFormat = "@"
$ws.Range("J4").Value = ""
$ws.Range("J4").Style = "Normal"
$ws.Range("P4").Value = "S om järnvägen - 2, Vg"
$ws.Range("Q4").Value = 432083.280685614
$ws.Range("R4").Value = 6419676.539718015

# Row 5
$ws.Range("A5").Value = 111482980
$ws.Range("I5").NumberFormat = "@"
$ws.Range("I5").Value = "10"
$ws.Range("I5").Style = "Normal"
$ws.Range("P5").Value = "S om järnvägen - 6, Vg"
$ws.Range("Q5").Value = 432048.2263952638
$ws.Range("R5").Value = 6419681.385014677

# Row 6
$ws.Range("A6").Value = 111482955
$ws.Range("I6").NumberFormat = "@"
$ws.Range("I6").Value = "70"
$ws.Range("I6").Style = "Normal"
$ws.Range("P6").Value = "S om järnvägen - 5, Vg"
$ws.Range("Q6").Value = 432064.1298546481
$ws.Range("R6").Value = 6419677.395781181

# Row 7
$ws.Range("A7").Value = 111482936
$ws.Range("I7").NumberFormat = "@"
$ws.Range("I7").Value = "25"
$ws.Range("I7").Style = "Normal"
$ws.Range("P7").Value = "S om järnvägen - 4, Vg"
$ws.Range("Q7").Value = 432073.5656663703
$ws.Range("R7").Value = 6419668.734013095

# Row 8
$ws.Range("A8").Value = 111483462
$ws.Range("I8").NumberFormat = "@"
$ws.Range("I8").Value = "45"
$ws.Range("I8").Style = "Normal"
$ws.Range("P8").Value = "S om järnvägen - 16, Vg"
$ws.Range("Q8").Value = 431654.0242198514
$ws.Range("R8").Value = 6419791.70470859

# Row 9
$ws.Range("A9").Value = 111490843
$ws.Range("B9").Value = 96348
$ws.Range("D9").Value = "VU"
$ws.Range("E9").Value = 220787
$ws.Range("F9").Value = "Knärot"
$ws.Range("G9").Value = "Goodyera repens"
$ws.Range("H9").Value = "(L.) R. Br."
$ws.Range("I9").NumberFormat = "@"
$ws.Range("I9").Value = "50"
$ws.Range("I9").Style = "Normal"
$ws.Range("J9").Value = "stjälkar/strån/skott"
$ws.Range("K9").Value = "fullt utvecklade blad"
$ws.Range("L9").NumberFormat = "@"
$ws.Range("L9").Value = ""
$ws.Range("L9").Style = "Normal"
$ws.Range("P9").Value = "S om järnvägen - 17, Vg"
$ws.Range("Q9").Value = 431803.2980747336
$ws.Range("R9").Value = 6419679.170503675
$ws.Range("AJ9").Value = $null
$ws.Range("AK9").Value = $null
$ws.Range("AM9").Value = $null
$ws.Range("AO9").Value = $null

# Row 10
$ws.Range("A10").Value = 111483105
$ws.Range("B10").Value = 73689
$ws.Range("D10").Value = "NT"
$ws.Range("E10").Value = 308
$ws.Range("F10").Value = "Brunpudrad nållav"
$ws.Range("G10").Value = "Chaenotheca gracillima"
$ws.Range("H10").Value = "(Vain.) Tibell"

# Row 11
$ws.Range("A11").Value = 111491187
$ws.Range("I11").NumberFormat = "@"
$ws.Range("I11").Value = "60"
$ws.Range("I11").Style = "Normal"
$ws.Range("K11").Value = "blomning"
$ws.Range("P11").Value = "S om järnvägen - 18, Vg"
$ws.Range("Q11").Value = 431829.514510141
$ws.Range("R11").Value = 6419749.394753682

# Row 12
$ws.Range("A12").Value = 111483037
$ws.Range("B12").Value = 96348
$ws.Range("D12").Value = "VU"
$ws.Range("E12").Value = 220787
$ws.Range("F12").Value = "Knärot"
$ws.Range("G12").Value = "Goodyera repens"
$ws.Range("H12").Value = "(L.) R. Br."
$ws.Range("I12").NumberFormat = "@"
$ws.Range("I12").Value = "60"
$ws.Range("I12").Style = "Normal"
$ws.Range("J12").Value = "stjälkar/strån/skott"
$ws.Range("K12").Value = "blomning"
$ws.Range("L12").NumberFormat = "@"
$ws.Range("L12").Value = ""
$ws.Range("L12").Style = "Normal"
$ws.Range("P12").Value = "S om järnvägen - 7, Vg"
$ws.Range("Q12").Value = 432060.6482816387
$ws.Range("R12").Value = 6419660.45125766
$ws.Range("AM12").Value = $null
$ws.Range("AO12").Value = $null

# Row 13
$ws.Range("A13").Value = 111483140
$ws.Range("B13").Value = 73683
$ws.Range("D13").Value = "LC"
$ws.Range("E13").Value = 306
$ws.Range("F13").Value = "Kornig nållav"
$ws.Range("G13").Value = "Chaenotheca chlorella"
$ws.Range("H13").Value = "(Ach.) Müll.Arg."
$ws.Range("I13").NumberFormat = "@"
$ws.Range("I13").Value = ""
$ws.Range("I13").Style = "Normal"
$ws.Range("J13").NumberFormat = "@"
$ws.Range("J13").Value = ""
$ws.Range("J13").Style = "Normal"
$ws.Range("K13").NumberFormat = "@"
$ws.Range("K13").Value = ""
$ws.Range("K13").Style = "Normal"
$ws.Range("L13").Value = $null
$ws.Range("P13").Value = "S om järnvägen - 9, Vg"
$ws.Range("Q13").Value = 431942.9372677525
$ws.Range("R13").Value = 6419625.784949708
$ws.Range("AJ13").Value = "tall"
$ws.Range("AK13").Value = "Pinus sylvestris"
$ws.Range("AM13").Value = "Stående död trädstam/högstubbe"
$ws.Range("AO13").Value = "Standing dead tree/snags # Pinus sylvestris"

# Row 14
$ws.Range("A14").Value = 111491635
$ws.Range("I14").NumberFormat = "@"
$ws.Range("I14").Value = "10"
$ws.Range("I14").Style = "Normal"
$ws.Range("P14").Value = "S om järnvägen - 21, Vg"
$ws.Range("Q14").Value = 431859.6228004749
$ws.Range("R14").Value = 6419672.898494411

# Row 15
$ws.Range("A15").Value = 111483381
$ws.Range("B15").Value = 73689
$ws.Range("D15").Value = "NT"
$ws.Range("E15").Value = 308
$ws.Range("F15").Value = "Brunpudrad nållav"
$ws.Range("G15").Value = "Chaenotheca gracillima"
$ws.Range("H15").Value = "(Vain.) Tibell"
$ws.Range("I15").NumberFormat = "@"
$ws.Range("I15").Value = ""
$ws.Range("I15").Style = "Normal"
$ws.Range("J15").NumberFormat = "@"
$ws.Range("J15").Value = ""
$ws.Range("J15").Style = "Normal"
$ws.Range("K15").NumberFormat = "@"
$ws.Range("K15").Value = ""
$ws.Range("K15").Style = "Normal"
$ws.Range("L15").Value = $null
$ws.Range("P15").Value = "S om järnvägen - 14, Vg"
$ws.Range("Q15").Value = 431754.10213514
$ws.Range("R15").Value = 6419728.893211351
$ws.Range("AJ15").Value = "tall"
$ws.Range("AK15").Value = "Pinus sylvestris"
$ws.Range("AM15").Value = "Stående död trädstam/högstubbe"
$ws.Range("AO15").Value = "Standing dead tree/snags # Pinus sylvestris"

# Row 16
$ws.Range("A16").Value = 111483197
$ws.Range("P16").Value = "S om järnvägen - 11, Vg"
$ws.Range("Q16").Value = 431937.082796899
$ws.Range("R16").Value = 6419625.884406033
$ws.Range("AJ16").Value = $null
$ws.Range("AK16").Value = $null
$ws.Range("AO16").Value = "Standing dead tree/snags"

# Row 18
$ws.Range("A18").Value = 111483300
$ws.Range("P18").Value = "S om järnvägen - 12, Vg"
$ws.Range("Q18").Value = 431888.091041417
$ws.Range("R18").Value = 6419625.122914318

# Row 19
$ws.Range("A19").Value = 111483107
$ws.Range("B19").Value = 73681
$ws.Range("D19").Value = "LC"
$ws.Range("E19").Value = 6439
$ws.Range("F19").Value = "Gulnål"
$ws.Range("G19").Value = "Chaenotheca brachypoda"
$ws.Range("H19").Value = "(Ach.) Tibell"
$ws.Range("P19").Value = "S om järnvägen - 8, Vg"
$ws.Range("Q19").Value = 431947.1499479365
$ws.Range("R19").Value = 6419623.056550305

# Row 20
$ws.Range("A20").Value = 111661838
$ws.Range("Q20").Value = 431799.2483237319
$ws.Range("R20").Value = 6419691.460736625

# Row 21
$ws.Range("A21").Value = 111661840
$ws.Range("B21").Value = 89793
$ws.Range("D21").Value = "LC"
$ws.Range("E21").Value = 4217
$ws.Range("F21").Value = "Blodticka"
$ws.Range("G21").Value = "Meruliopsis taxicola"
$ws.Range("H21").Value = "(Pers.:Fr.) Bondartsev"
$ws.Range("J21").NumberFormat = "@"
$ws.Range("J21").Value = ""
$ws.Range("J21").Style = "Normal"
$ws.Range("K21").NumberFormat = "@"
$ws.Range("K21").Value = ""
$ws.Range("K21").Style = "Normal"
$ws.Range("N21").NumberFormat = "@"
$ws.Range("N21").Value = ""
$ws.Range("N21").Style = "Normal"
$ws.Range("Q21").Value = 431768.994999051
$ws.Range("R21").Value = 6419728.1081824
$ws.Range("AF21").NumberFormat = "@"
$ws.Range("AF21").Value = ""
$ws.Range("AF21").Style = "Normal"
$ws.Range("AJ21").Value = "tall"
$ws.Range("AK21").Value = "Pinus sylvestris"
$ws.Range("AO21").Value = "Pinus sylvestris"

# Row 22
$ws.Range("A22").Value = 111661832
$ws.Range("Q22").Value = 432076.4609239195
$ws.Range("R22").Value = 6419682.500295377

# Row 23
$ws.Range("A23").Value = 111661831
$ws.Range("B23").Value = 96348
$ws.Range("D23").Value = "VU"
$ws.Range("E23").Value = 220787
$ws.Range("F23").Value = "Knärot"
$ws.Range("G23").Value = "Goodyera repens"
$ws.Range("H23").Value = "(L.) R. Br."
$ws.Range("J23").Value = $null
$ws.Range("K23").Value = $null
$ws.Range("N23").Value = $null
$ws.Range("Q23").Value = 432080.3854477856
$ws.Range("R23").Value = 6419662.773410858
$ws.Range("AF23").Value = $null
$ws.Range("AJ23").Value = $null
$ws.Range("AK23").Value = $null
$ws.Range("AO23").Value = $null
